$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the default (unstyled) cell style used by the data cells so that
# forcing text NumberFormat on numeric-looking values does not change the
# persisted cell style index.
$normalStyle = $ws.Range("D2").Style

$ws.Range("D2").Value = "26.150.41"
$ws.Range("E2").Value = "  +3.67%  "

$ws.Range("D3").Value = "1.602.88"
$ws.Range("E3").Value = "  +3.48%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.50"
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = "  +2.68%  "

$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.486"
$ws.Range("D7").Style = $normalStyle
$ws.Range("E7").Value = "  +2.14%  "

$ws.Range("E8").Value = "  +2.57%  "

$ws.Range("E9").Value = "  +1.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.98"
$ws.Range("D10").Style = $normalStyle
$ws.Range("E10").Value = "  +1.29%  "

$ws.Range("E11").Value = "  +4.91%  "

$ws.Range("D12").Value = "1.827.44"

$ws.Range("D13").Value = "1.602.57"
$ws.Range("E13").Value = "  +3.70%  "

$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.510"
$ws.Range("D15").Style = $normalStyle
$ws.Range("E15").Value = "  +1.48%  "

$ws.Range("D16").Value = "26.135.91"
$ws.Range("E16").Value = "  +3.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.45"
$ws.Range("D17").Style = $normalStyle
$ws.Range("E17").Value = "  +3.33%  "

$ws.Range("D18").Value = "0.0₃0721"
$ws.Range("E18").Value = "  +2.10%  "

$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "204.39"
$ws.Range("D20").Style = $normalStyle
$ws.Range("E20").Value = "  +10.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.23"
$ws.Range("D21").Style = $normalStyle
$ws.Range("E21").Value = "  +3.40%  "

$ws.Range("E22").Value = "  +0.76%  "

$ws.Range("E23").Value = "  +2.79%  "

$ws.Range("E24").Value = "  +11.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.59"
$ws.Range("D25").Style = $normalStyle
$ws.Range("E25").Value = "  +1.78%  "

$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.123"
$ws.Range("D27").Style = $normalStyle
$ws.Range("E27").Value = "  -3.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.18"
$ws.Range("D28").Style = $normalStyle
$ws.Range("E28").Value = "  +2.74%  "

$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("E30").Value = "  +1.59%  "

$ws.Range("E31").Value = "  +1.98%  "

$ws.Range("E32").Value = "  +3.33%  "

$ws.Range("E33").Value = "  +0.46%  "

$ws.Range("E34").Value = "  +1.63%  "

$ws.Range("E35").Value = "  +1.90%  "

$ws.Range("D36").Value = "1.118.74"
$ws.Range("E36").Value = "  +3.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0163"
$ws.Range("D37").Style = $normalStyle
$ws.Range("E37").Value = "  +9.72%  "

$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("E39").Value = "  +3.33%  "

$ws.Range("E40").Value = "  +2.54%  "

$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.781"
$ws.Range("D42").Style = $normalStyle
$ws.Range("E42").Value = "  -2.11%  "

$ws.Range("D43").Value = "1.739.22"
$ws.Range("E43").Value = "  +3.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.12"
$ws.Range("D44").Style = $normalStyle
$ws.Range("E44").Value = "  +1.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.92"
$ws.Range("D45").Style = $normalStyle

$ws.Range("E46").Value = "  +3.96%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.47"
$ws.Range("D47").Style = $normalStyle

$ws.Range("E48").Value = "  +0.48%  "

$ws.Range("E49").Value = "  +1.31%  "

$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").Value = "0.0₇0924"
$ws.Range("E51").Value = "  -16.43%  "
